$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9996787309646606
$ws.Range("B1").Value = 4.185778141021729
$ws.Range("C1").Value = 2.407150506973267
$ws.Range("D1").Value = 1.750830888748169
$ws.Range("E1").Value = 1.375241041183472
